# Two rows were removed from the missing-data worksheet ("RM 232" and
# "SC 92"), which shifts every subsequent row up. After the shift, a few
# of the previously-blank/filled cells differ from a straight shift:
#   - "SC 5"   (now row 26) column C gains a value (10.8)
#   - "SC 101" (now row 27) column C becomes blank
#   - "SC 232" (now row 33) column D gains a value (-14.1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "RM 232" row (row 26).
$ws.Rows.Item(26).Delete()

# After the above shift, the "SC 92" row is now row 27; delete it too.
$ws.Rows.Item(27).Delete()

# Fill in / clear the remaining cells that changed value (not just position).
$ws.Range("C26").Value = 10.8
$ws.Range("C27").Value = ""
$ws.Range("D33").Value = -14.1
